# Natmi following Dr Hou advice:
# Regenerate the Dhh-Boc sending/receiving cluster table to include the new
# "ECs" sending cluster alongside the existing "FAPs" and "sCs" clusters.
# The sheet becomes a full 3 (senders) x 3 (targets) cross join (rows 2-10)
# of ECs/FAPs/sCs, with ligand "Dhh" and receptor "Boc" fixed, and refreshed
# statistics for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write string columns (A-D) in column-major order (all of column A, then all
# of column B, etc.) so that any newly-introduced text is (re)discovered in
# the same order as the source data: ECs, FAPs, sCs, Dhh, Boc.
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "sCs"
$ws.Range("A9").Value = "sCs"
$ws.Range("A10").Value = "sCs"

$ws.Range("B2").Value = "Dhh"
$ws.Range("B3").Value = "Dhh"
$ws.Range("B4").Value = "Dhh"
$ws.Range("B5").Value = "Dhh"
$ws.Range("B6").Value = "Dhh"
$ws.Range("B7").Value = "Dhh"
$ws.Range("B8").Value = "Dhh"
$ws.Range("B9").Value = "Dhh"
$ws.Range("B10").Value = "Dhh"

$ws.Range("C2").Value = "Boc"
$ws.Range("C3").Value = "Boc"
$ws.Range("C4").Value = "Boc"
$ws.Range("C5").Value = "Boc"
$ws.Range("C6").Value = "Boc"
$ws.Range("C7").Value = "Boc"
$ws.Range("C8").Value = "Boc"
$ws.Range("C9").Value = "Boc"
$ws.Range("C10").Value = "Boc"

$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "sCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "sCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("D10").Value = "sCs"

# Write numeric columns (E-T)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.179395
$ws.Range("H2").Value = 6.538185
$ws.Range("I2").Value = 0.4845018986408914
$ws.Range("J2").Value = 0.4845018986408914
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.412261333333333
$ws.Range("N2").Value = 4.236784
$ws.Range("O2").Value = 0.0207702079890636
$ws.Range("P2").Value = 0.0207702079890636
$ws.Range("Q2").Value = 3.07787528856
$ws.Range("R2").Value = 27.70087759704
$ws.Range("S2").Value = 0.01006320520586752
$ws.Range("T2").Value = 0.01006320520586752

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.179395
$ws.Range("H3").Value = 6.538185
$ws.Range("I3").Value = 0.4845018986408914
$ws.Range("J3").Value = 0.4845018986408914
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 60.17634833333333
$ws.Range("N3").Value = 180.529045
$ws.Range("O3").Value = 0.8850169875823316
$ws.Range("P3").Value = 0.8850169875823317
$ws.Range("Q3").Value = 131.148032675925
$ws.Range("R3").Value = 1180.332294083325
$ws.Range("S3").Value = 0.4287924108130818
$ws.Range("T3").Value = 0.4287924108130819

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 2.179395
$ws.Range("H4").Value = 6.538185
$ws.Range("I4").Value = 0.4845018986408914
$ws.Range("J4").Value = 0.4845018986408914
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.405959
$ws.Range("N4").Value = 19.217877
$ws.Range("O4").Value = 0.0942128044286047
$ws.Range("P4").Value = 0.09421280442860472
$ws.Range("Q4").Value = 13.961115014805
$ws.Range("R4").Value = 125.650035133245
$ws.Range("S4").Value = 0.04564628262194196
$ws.Range("T4").Value = 0.04564628262194197

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7348883333333333
$ws.Range("H5").Value = 2.204665
$ws.Range("I5").Value = 0.1633732264179005
$ws.Range("J5").Value = 0.1633732264179005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.412261333333333
$ws.Range("N5").Value = 4.236784
$ws.Range("O5").Value = 0.0207702079890636
$ws.Range("P5").Value = 0.0207702079890636
$ws.Range("Q5").Value = 1.037854377484444
$ws.Range("R5").Value = 9.34068939736
$ws.Range("S5").Value = 0.003393295892544172
$ws.Range("T5").Value = 0.003393295892544173

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7348883333333333
$ws.Range("H6").Value = 2.204665
$ws.Range("I6").Value = 0.1633732264179005
$ws.Range("J6").Value = 0.1633732264179005
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 60.17634833333333
$ws.Range("N6").Value = 180.529045
$ws.Range("O6").Value = 0.8850169875823316
$ws.Range("P6").Value = 0.8850169875823317
$ws.Range("Q6").Value = 44.22289633276944
$ws.Range("R6").Value = 398.006066994925
$ws.Range("S6").Value = 0.1445880806959765
$ws.Range("T6").Value = 0.1445880806959765

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7348883333333333
$ws.Range("H7").Value = 2.204665
$ws.Range("I7").Value = 0.1633732264179005
$ws.Range("J7").Value = 0.1633732264179005
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.405959
$ws.Range("N7").Value = 19.217877
$ws.Range("O7").Value = 0.0942128044286047
$ws.Range("P7").Value = 0.09421280442860472
$ws.Range("Q7").Value = 4.707664532911666
$ws.Range("R7").Value = 42.368980796205
$ws.Range("S7").Value = 0.01539184982937981
$ws.Range("T7").Value = 0.01539184982937982

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.583934333333333
$ws.Range("H8").Value = 4.751803
$ws.Range("I8").Value = 0.3521248749412083
$ws.Range("J8").Value = 0.3521248749412082
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.412261333333333
$ws.Range("N8").Value = 4.236784
$ws.Range("O8").Value = 0.0207702079890636
$ws.Range("P8").Value = 0.0207702079890636
$ws.Range("Q8").Value = 2.236929213505778
$ws.Range("R8").Value = 20.132362921552
$ws.Range("S8").Value = 0.007313706890651903
$ws.Range("T8").Value = 0.007313706890651903

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.583934333333333
$ws.Range("H9").Value = 4.751803
$ws.Range("I9").Value = 0.3521248749412083
$ws.Range("J9").Value = 0.3521248749412082
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 60.17634833333333
$ws.Range("N9").Value = 180.529045
$ws.Range("O9").Value = 0.8850169875823316
$ws.Range("P9").Value = 0.8850169875823317
$ws.Range("Q9").Value = 95.31538417979277
$ws.Range("R9").Value = 857.838457618135
$ws.Range("S9").Value = 0.3116364960732734
$ws.Range("T9").Value = 0.3116364960732734

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.583934333333333
$ws.Range("H10").Value = 4.751803
$ws.Range("I10").Value = 0.3521248749412083
$ws.Range("J10").Value = 0.3521248749412082
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.405959
$ws.Range("N10").Value = 19.217877
$ws.Range("O10").Value = 0.0942128044286047
$ws.Range("P10").Value = 0.09421280442860472
$ws.Range("Q10").Value = 10.14661839802567
$ws.Range("R10").Value = 91.31956558223101
$ws.Range("S10").Value = 0.03317467197728294
$ws.Range("T10").Value = 0.03317467197728294
